$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "67.604.74"
$ws.Cells.Item(2, 5).Value = "  +1.80%  "
$ws.Cells.Item(3, 4).Value = "2.618.33"
$ws.Cells.Item(4, 5).Value = "  -0.08%  "
$ws.Cells.Item(5, 4).Value = "602.41"
$ws.Cells.Item(5, 5).Value = "  +2.24%  "
$ws.Cells.Item(6, 4).Value = "153.97"
$ws.Cells.Item(6, 5).Value = "  +0.72%  "
$ws.Cells.Item(7, 5).Value = "  -0.01%  "
$ws.Cells.Item(8, 4).Value = "0.550"
$ws.Cells.Item(8, 5).Value = "  +1.99%  "
$ws.Cells.Item(9, 4).Value = "2.615.97"
$ws.Cells.Item(9, 5).Value = "  +1.61%  "
$ws.Cells.Item(10, 4).Value = "0.127"
$ws.Cells.Item(10, 5).Value = "  +12.12%  "
$ws.Cells.Item(11, 5).Value = "  +0.86%  "
$ws.Cells.Item(12, 5).Value = "  +1.51%  "
$ws.Cells.Item(13, 5).Value = "  +0.59%  "
$ws.Cells.Item(14, 4).Value = "27.98"
$ws.Cells.Item(14, 5).Value = "  +0.43%  "
$ws.Cells.Item(15, 4).Value = "0.0000187"
$ws.Cells.Item(15, 5).Value = "  +4.82%  "
$ws.Cells.Item(16, 4).Value = "3.096.15"
$ws.Cells.Item(16, 5).Value = "  +1.52%  "
$ws.Cells.Item(17, 4).Value = "67.661.23"
$ws.Cells.Item(17, 5).Value = "  +2.03%  "
$ws.Cells.Item(18, 4).Value = "2.620.42"
$ws.Cells.Item(18, 5).Value = "  +1.75%  "
$ws.Cells.Item(19, 4).Value = "11.29"
$ws.Cells.Item(19, 5).Value = "  +0.78%  "
$ws.Cells.Item(20, 4).Value = "363.63"
$ws.Cells.Item(20, 5).Value = "  +3.70%  "
$ws.Cells.Item(21, 4).Value = "7.65"
$ws.Cells.Item(21, 5).Value = "  -1.12%  "
$ws.Cells.Item(22, 5).Value = "  -0.17%  "
$ws.Cells.Item(23, 4).Value = "2.13"
$ws.Cells.Item(23, 5).Value = "  +6.18%  "
$ws.Cells.Item(24, 4).Value = "0.999"
$ws.Cells.Item(24, 5).Value = "  -0.09%  "
$ws.Cells.Item(25, 4).Value = "70.11"
$ws.Cells.Item(25, 5).Value = "  +4.33%  "
$ws.Cells.Item(26, 4).Value = "10.01"
$ws.Cells.Item(26, 5).Value = "  -2.66%  "
$ws.Cells.Item(27, 4).Value = "0.0000106"
$ws.Cells.Item(27, 5).Value = "  +4.04%  "
$ws.Cells.Item(28, 4).Value = "2.745.37"
$ws.Cells.Item(28, 5).Value = "  +1.40%  "
$ws.Cells.Item(29, 4).Value = "581.94"
$ws.Cells.Item(29, 5).Value = "  -1.48%  "
$ws.Cells.Item(30, 2).Value = "Fetch.AI"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(30, 4).Value = "1.44"
$ws.Cells.Item(30, 5).Value = "  +0.57%  "
$ws.Cells.Item(31, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(31, 4).Value = "7.95"
$ws.Cells.Item(31, 5).Value = "  +0.03%  "
$ws.Cells.Item(32, 2).Value = "PancakeSwap"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(32, 4).Value = "1.87"
$ws.Cells.Item(32, 5).Value = "  +1.23%  "
$ws.Cells.Item(33, 2).Value = "Kaspa"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(33, 4).Value = "0.131"
$ws.Cells.Item(33, 5).Value = "  -1.19%  "
$ws.Cells.Item(34, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(34, 4).Value = "0.999"
$ws.Cells.Item(34, 5).Value = "  +0.03%  "
$ws.Cells.Item(35, 2).Value = "ImmutableX"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(35, 4).Value = "1.54"
$ws.Cells.Item(35, 5).Value = "  -0.86%  "
$ws.Cells.Item(36, 2).Value = "NEARProtocol"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(36, 4).Value = "4.98"
$ws.Cells.Item(36, 5).Value = "  +0.43%  "
$ws.Cells.Item(37, 2).Value = "Monero"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(37, 4).Value = "157.28"
$ws.Cells.Item(37, 5).Value = "  +2.75%  "
$ws.Cells.Item(38, 2).Value = "EthereumClassic"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(38, 4).Value = "19.45"
$ws.Cells.Item(38, 5).Value = "  +1.95%  "
$ws.Cells.Item(39, 2).Value = "PolygonEcosystemToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Cells.Item(39, 4).Value = "0.372"
$ws.Cells.Item(39, 5).Value = "  +1.38%  "
$ws.Cells.Item(40, 2).Value = "RenderToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Cells.Item(40, 4).Value = "5.40"
$ws.Cells.Item(40, 5).Value = "  +0.28%  "
$ws.Cells.Item(41, 2).Value = "Stacks"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(41, 4).Value = "1.85"
$ws.Cells.Item(41, 5).Value = "  +4.33%  "
$ws.Cells.Item(42, 2).Value = "dogwifhat"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(42, 4).Value = "2.70"
$ws.Cells.Item(42, 5).Value = "  +4.88%  "
$ws.Cells.Item(43, 2).Value = "OKB"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(43, 4).Value = "41.12"
$ws.Cells.Item(43, 5).Value = "  -0.30%  "
$ws.Cells.Item(44, 2).Value = "USDe"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(44, 4).Value = "1.00"
$ws.Cells.Item(44, 5).Value = "  +0.03%  "
$ws.Cells.Item(45, 4).Value = "16.41"
$ws.Cells.Item(45, 5).Value = "  +0.05%  "
$ws.Cells.Item(46, 2).Value = "Aave"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(46, 4).Value = "156.90"
$ws.Cells.Item(46, 5).Value = "  +1.31%  "
$ws.Cells.Item(47, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(47, 4).Value = "0.0₆0292"
$ws.Cells.Item(47, 5).Value = "  -4.16%  "
$ws.Cells.Item(48, 2).Value = "Filecoin"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(48, 4).Value = "3.77"
$ws.Cells.Item(48, 5).Value = "  +1.04%  "
$ws.Cells.Item(49, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(49, 4).Value = "21.05"
$ws.Cells.Item(49, 5).Value = "  +0.72%  "
$ws.Cells.Item(50, 2).Value = "Mantle"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(50, 4).Value = "0.623"
$ws.Cells.Item(50, 5).Value = "  +2.01%  "
$ws.Cells.Item(51, 2).Value = "Hedera"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(51, 4).Value = "0.0540"
$ws.Cells.Item(51, 5).Value = "  -1.79%  "
